$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp + per-country data rows to reflect the 13 Abril 2020 17:52 refresh.
# Row layout: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$ws.Range("A1").Value = 'Datos actualizados a 13 de Abril de 2020 a las 17:52'

# Row 4: Estados Unidos
$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 564171
$ws.Range("C4").Value = 3871
$ws.Range("D4").Value = 33728
$ws.Range("E4").Value = 507598
$ws.Range("F4").Value = 11807
$ws.Range("G4").Value = 740
$ws.Range("H4").Value = 22845

# Row 17: Brasil
$ws.Range("A17").Value = 'Brasil'
$ws.Range("B17").Value = 22720
$ws.Range("C17").Value = 528
$ws.Range("D17").Value = 173
$ws.Range("E17").Value = 21278
$ws.Range("F17").Value = 296
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 1269

# Row 21: Israel
$ws.Range("A21").Value = 'Israel'
$ws.Range("B21").Value = 11235
$ws.Range("C21").Value = 90
$ws.Range("D21").Value = 1689
$ws.Range("E21").Value = 9433
$ws.Range("F21").Value = 181
$ws.Range("G21").Value = 10
$ws.Range("H21").Value = 113

# Row 25: India
$ws.Range("A25").Value = 'India'
$ws.Range("B25").Value = 9635
$ws.Range("C25").Value = 430
$ws.Range("D25").Value = 1096
$ws.Range("E25").Value = 8208
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 331

# Row 32: Noruega
$ws.Range("A32").Value = 'Noruega'
$ws.Range("B32").Value = 6547
$ws.Range("C32").Value = 22
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 6382
$ws.Range("F32").Value = 59
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 133

# Row 35: Chequia
$ws.Range("A35").Value = 'Chequia'
$ws.Range("B35").Value = 6022
$ws.Range("C35").Value = 31
$ws.Range("D35").Value = 519
$ws.Range("E35").Value = 5360
$ws.Range("F35").Value = 87
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 143

# Row 82: Cuba
$ws.Range("A82").Value = 'Cuba'
$ws.Range("B82").Value = 726
$ws.Range("C82").Value = 57
$ws.Range("D82").Value = 121
$ws.Range("E82").Value = 584
$ws.Range("F82").Value = 11
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 21

# Row 83: Crucero
$ws.Range("A83").Value = 'Crucero'
$ws.Range("B83").Value = 712
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 619
$ws.Range("E83").Value = 82
$ws.Range("F83").Value = 10
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 11

# Row 84: Tunez
$ws.Range("A84").Value = 'Tunez'
$ws.Range("B84").Value = 707
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 43
$ws.Range("E84").Value = 633
$ws.Range("F84").Value = 89
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 31

# Row 85: Bulgaria
$ws.Range("A85").Value = 'Bulgaria'
$ws.Range("B85").Value = 685
$ws.Range("C85").Value = 10
$ws.Range("D85").Value = 71
$ws.Range("E85").Value = 582
$ws.Range("F85").Value = 36
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 32

# Row 87: Republica de Chipre
$ws.Range("A87").Value = 'Republica de Chipre'
$ws.Range("B87").Value = 662
$ws.Range("C87").Value = 29
$ws.Range("D87").Value = 65
$ws.Range("E87").Value = 585
$ws.Range("F87").Value = 8
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 12

# Row 108: Estado de Palestina
$ws.Range("A108").Value = 'Estado de Palestina'
$ws.Range("B108").Value = 308
$ws.Range("C108").Value = 18
$ws.Range("D108").Value = 58
$ws.Range("E108").Value = 248
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 2

# Row 109: Senegal
$ws.Range("A109").Value = 'Senegal'
$ws.Range("B109").Value = 291
$ws.Range("C109").Value = 11
$ws.Range("D109").Value = 178
$ws.Range("E109").Value = 111
$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 2

# Row 130: Mali
$ws.Range("A130").Value = 'Mali'
$ws.Range("B130").Value = 123
$ws.Range("C130").Value = 18
$ws.Range("D130").Value = 26
$ws.Range("E130").Value = 87
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 10

# Row 131: Camboya
$ws.Range("A131").Value = 'Camboya'
$ws.Range("B131").Value = 122
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 77
$ws.Range("E131").Value = 45
$ws.Range("F131").Value = 1
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 0

# Row 132: Trinidad yTobago
$ws.Range("A132").Value = 'Trinidad yTobago'
$ws.Range("B132").Value = 113
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 16
$ws.Range("E132").Value = 89
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 8

# Row 133: Madagascar
$ws.Range("A133").Value = 'Madagascar'
$ws.Range("B133").Value = 106
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 21
$ws.Range("E133").Value = 85
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

# Row 143: Liberia
$ws.Range("A143").Value = 'Liberia'
$ws.Range("B143").Value = 59
$ws.Range("C143").Value = 9
$ws.Range("D143").Value = 4
$ws.Range("E143").Value = 49
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 6

# Row 144: Gabon
$ws.Range("A144").Value = 'Gabon'
$ws.Range("B144").Value = 57
$ws.Range("C144").Value = 8
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 55
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 1

# Row 145: Bermudas
$ws.Range("A145").Value = 'Bermudas'
$ws.Range("B145").Value = 57
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 29
$ws.Range("E145").Value = 24
$ws.Range("F145").Value = 2
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 4

# Row 146: Uganda
$ws.Range("A146").Value = 'Uganda'
$ws.Range("B146").Value = 54
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 7
$ws.Range("E146").Value = 47
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# Row 147: Polinesia Francesa
$ws.Range("A147").Value = 'Polinesia Francesa'
$ws.Range("B147").Value = 53
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 0
$ws.Range("E147").Value = 53
$ws.Range("F147").Value = 1
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 0

# Row 148: Islas Caimanes
$ws.Range("A148").Value = 'Islas Caimanes'
$ws.Range("B148").Value = 53
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 6
$ws.Range("E148").Value = 46
$ws.Range("F148").Value = 3
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 1

# Row 208: Santo Tome y Principe
$ws.Range("A208").Value = 'Santo Tome y Principe'
$ws.Range("B208").Value = 4
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 0
$ws.Range("E208").Value = 4
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209: Sudan del Sur
$ws.Range("A209").Value = 'Sudan del Sur'
$ws.Range("B209").Value = 4
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 4
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0
